$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update "VALOR MORA" value (E11): 113880 -> 56940
$ws.Range("E11").Value = 56940

# 2. Update "Cant. Periodos" value (F13): 2 -> 1
$ws.Range("F13").Value = 1

# 3. Update the surviving period row's value from 2507 -> 2509 (row 16, column E)
$ws.Range("E16").Value = "2509"

# 4. Delete the now-obsolete second period row (row 17) entirely, shifting rows below up
$ws.Rows("17").Delete()
